$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Fix "conceptPath" values in column J: replace "/root/" with "/tbi/"
#    (commit: "always check with lowercase" / path prefix fix root -> tbi)
# -----------------------------------------------------------------
$conceptRange = $ws.Range("J1:J159")
$conceptRange.Replace("/root/", "/tbi/") | Out-Null

# Special case: row 6 originally had "root/creactive/glucose" (missing the
# leading slash), which becomes "/tbi/creactive/glucose".
$ws.Range("J6").Value = "/tbi/creactive/glucose"

# -----------------------------------------------------------------
# 2) Populate the "methodology" column (K) with "mip-cde" for every data
#    row that was still missing it (rows 25 through 159).
# -----------------------------------------------------------------
for ($r = 25; $r -le 159; $r++) {
    $ws.Range("K$r").Value = "mip-cde"
}

# Match the look of the existing methodology cells (K2:K24), which wrap text.
$ws.Range("K25:K159").WrapText = $true
